$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2019")

$values = @(
    "sample id",
    "1-58_170m_NO2_ctl",
    "1-58_170m_NO2_+P",
    "1-58_570m_NO2_ctl",
    "1_58_570m_NO2_+P",
    "1-58_570m_NO2_top",
    "1-58_370m_ctl",
    "1-58_370m_+P",
    "1-58_370m_NO2_top",
    "1-58_102m_ctl",
    "1-58_102m_+P",
    "1-58_102m_top",
    "1-58_770m_net_+P",
    "1-58_770m_net_top",
    "1-58_770m_net_nw",
    "4-59_130m_NO2_+P",
    "4-59_130m_NO2_ctl",
    "4-59_130m_NO2_top",
    "4-59_300m_Fe_top_split",
    "4-59_700m_NO2_+P",
    "4-59_700m_NO2_top",
    "4-59_700m_NO2_ctl",
    "4-59_965m_top",
    "4-59_965m_+P",
    "4-59_965m_ctl",
    "4-52_173m_+P",
    "4-52_173m_ctl",
    "4-52_173m_top",
    "4-52_965m_top/+P",
    "3-54_402m_ctl",
    "3-54_402m_top",
    "3-54_171m_+P",
    "3-54_200m_ctl",
    "3-54_200m_top",
    "1-50_100m_NO2_ctl",
    "1-50_300m_NO2_top",
    "1-50_300m_NO2_+P",
    "1-50_100m_NO2_top",
    "1-50_100m_NO2_+P",
    "1-50_300m_NO2_ctl",
    "4-55_490m_top/+P",
    "4-55_390m_top_Clint",
    "4-55_90m_+P",
    "4-55_490m_ctl",
    "4-55_90m_ctl",
    "3-51_520m_NO2_+P",
    "3-51_520m_NO2_ctl",
    "1-53_278m_+P",
    "1-53_107m_+P",
    "1-53_107m_top",
    "1-53_107m_ctl"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

$ws.Range("B11").Select()

Write-Output "applied"
